# Change "...студенту Федорову Льву Александровичу" into
# "...студентам Федорову Льву Александровичу и Садику Назару Самировичу"
#
# Concretely (per the target OOXML diff) the paragraph's runs become:
#   "с" | "тудент" | "ам" | " Федорову Льву Александровичу" | " и Садику Назару Самировичу"
# (all five runs share the same formatting: color 000000).

$d = $word.ActiveDocument

# 1) Locate the run whose text is "туденту".
$findRng = $d.Content
$found = $findRng.Find.Execute("туденту", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'туденту' in the document"
}
$wordStart = $findRng.Start
$wordEnd = $findRng.End

# Find the paragraph that contains it, so its end can be tracked live as we edit.
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $wordStart -and $p.Range.End -ge $wordEnd) {
        $para = $p
        break
    }
}
if ($null -eq $para) {
    throw "Could not locate the paragraph containing 'туденту'"
}

# 2) Shrink "туденту" -> "тудент" by dropping the trailing "у".
$dropU = $d.Range($wordEnd - 1, $wordEnd)
$dropU.Text = ""
$afterTudent = $wordEnd - 1

# 3) Insert "ам" right after "тудент" (=> "студентам").
$insAm = $d.Range($afterTudent, $afterTudent)
$insAm.InsertAfter("ам")
$afterAm = $afterTudent + 2

# The text " Федорову Льву Александровичу" that used to follow "туденту"
# is still in place right after "ам"; remember where it ends (end of the
# paragraph's text, i.e. just before the paragraph mark).
$afterFedorov = $para.Range.End - 1

# 4) Append the new tail " и Садику Назару Самировичу" at the very end of
#    the paragraph.
$tail = $d.Range($afterFedorov, $afterFedorov)
$tail.InsertAfter(" и Садику Назару Самировичу")

# 5) The edits above cause Word to coalesce same-formatted runs in this
#    paragraph into a single run. Re-split it back into the five pieces
#    described above by toggling Bold off/on across each final segment's
#    boundaries (this forces a run split without any visible formatting
#    change, since Bold ends up false again on every segment).
$paraStart = $wordStart - 1
$segments = @(
    @($paraStart, $wordStart),     # "с"
    @($wordStart, $afterTudent),   # "тудент"
    @($afterTudent, $afterAm),     # "ам"
    @($afterAm, $afterFedorov),    # " Федорову Льву Александровичу"
    @($afterFedorov, ($para.Range.End - 1))  # " и Садику Назару Самировичу"
)
foreach ($seg in $segments) {
    $r = $d.Range($seg[0], $seg[1])
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

Write-Output "Done. Paragraph now reads:"
Write-Output $para.Range.Text
